$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (serial 45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update price values in column D
$ws.Range("D22").Value = 12264
$ws.Range("D23").Value = 13894
$ws.Range("D24").Value = 18098
$ws.Range("D25").Value = 20048
